# Updating RTM to reflect modified requirements according to review.
# - CYRS/SRS/HSI status columns move from "Released"/"Draft" to "Proposed"
#   for every existing requirement row.
# - A new HSI Requirment ID is recorded for each existing CRS row.
# - A new CYRS requirement (REQ_PO2EBL_CYRS_03_V1.3) and its related HSI
#   requirement (REQ_PO2EBL_HSI_05_V01) are appended below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Req_PO2EBL_CRS_01_V01) ---------------------------------------
$ws.Range("D2").Value = "Proposed"
$ws.Range("G2").Value = "Proposed"
$ws.Range("J2").Value = "Proposed"
$ws.Range("H2").Value = "REQ_PO2EBL_HSI_01_V01"

# --- Row 3 ------------------------------------------------------------
$ws.Range("D3").Value = "Proposed"
$ws.Range("G3").Value = "Proposed"
$ws.Range("J3").Value = "Proposed"
$ws.Range("H3").Value = "REQ_PO2EBL_HSI_02_V01"

# --- Row 4 ------------------------------------------------------------
$ws.Range("D4").Value = "Proposed"
$ws.Range("G4").Value = "Proposed"
$ws.Range("J4").Value = "Proposed"

# A blank spacer row is inserted above the old "Req_PO2EBL_CRS_02_V01" row,
# pushing it (and the separator row after it) one row down.
$ws.Rows(5).Insert()

# --- New row 9: additional CYRS/HSI requirement pair ---------------------
$ws.Range("H9").Value = "REQ_PO2EBL_HSI_05_V01"
$ws.Range("A9").Value = "REQ_PO2EBL_CYRS_03_V1.3"
$ws.Range("A9").Font.Name = "Calibri Light"

# --- Row 7 (previously row 6, "Req_PO2EBL_CRS_02_V01") -------------------
$ws.Range("D7").Value = "Proposed"
$ws.Range("G7").Value = "Proposed"
$ws.Range("J7").Value = "Proposed"
$ws.Range("H7").Value = "REQ_PO2EBL_HSI_03_V01"
$ws.Range("H7").Font.Name = "Times New Roman"
$ws.Range("H7").Font.Size = 10

# Finish the HSI id for row 4 (kept here so the shared-string order matches
# the order the requirements were actually entered in).
$ws.Range("H4").Value = "REQ_PO2EBL_HSI_04_V01"

# Leave the selection where the author left off editing.
$ws.Range("H7").Select() | Out-Null
